$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows appended to the bottom of the data (rows 11-13)
$newRows = @(
    @{ Row = 11; A = 9785.9500000000007; B = 9702.51;            C = 305.24; D = 307.86; E = $false; F = 0.86;  G = 42613.765451388892; H = $true  },
    @{ Row = 12; A = 9720.3799999999992; B = 9785.9500000000007;  C = 307.68; D = 305.63; E = $false; F = -0.67; G = 42614.672592592593; H = $false },
    @{ Row = 13; A = 9659.14;            B = 9720.3799999999992;  C = 307.96; D = 306.02; E = $false; F = -0.63; G = 42615.750023148146; H = $false }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.Value = $r.G
    $gCell.NumberFormat = "m/d/yy h:mm"

    $ws.Cells.Item($row, 8).Value = $r.H
}
